# Edit the SL fMRI experiment stimuli sheet (visual_run3_2.xlsx)
# Rebuilds the trial table (columns: image, trialnum, condition, word,
# location, repetition) for rows 2-49 with the new stimulus set
# (A.png..M.png letter images instead of AlienN.bmp) and removes the old
# trailing "blank"/condition-B filler rows, replacing them with real R
# trials so the table now runs trialnum 1-48 continuously.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header) is unchanged: image | trialnum | condition | word | location | repetition

$data = @(
    @("B.png", 1, "R", 0, 1, 1),
    @("A.png", 2, "R", 0, 2, 1),
    @("K.png", 3, "R", 0, 3, 1),
    @("D.png", 4, "R", 0, 1, 1),
    @("C.png", 5, "R", 0, 2, 1),
    @("M.png", 6, "R", 0, 3, 1),
    @("E.png", 7, "R", 0, 1, 1),
    @("C.png", 8, "R", 0, 2, 1),
    @("J.png", 9, "R", 0, 3, 1),
    @("F.png", 10, "R", 0, 1, 1),
    @("J.png", 11, "R", 0, 2, 1),
    @("L.png", 12, "R", 0, 3, 1),
    @("G.png", 13, "R", 0, 1, 1),
    @("G.png", 14, "R", 0, 2, 1),
    @("A.png", 15, "R", 0, 3, 1),
    @("H.png", 16, "R", 0, 1, 1),
    @("K.png", 17, "R", 0, 2, 1),
    @("M.png", 18, "R", 0, 3, 1),
    @("F.png", 19, "R", 0, 1, 1),
    @("K.png", 20, "R", 0, 2, 1),
    @("B.png", 21, "R", 0, 3, 1),
    @("M.png", 22, "R", 0, 1, 1),
    @("D.png", 23, "R", 0, 2, 1),
    @("B.png", 24, "R", 0, 3, 1),
    @("K.png", 25, "R", 0, 1, 1),
    @("L.png", 26, "R", 0, 2, 1),
    @("B.png", 27, "R", 0, 3, 1),
    @("G.png", 28, "R", 0, 1, 1),
    @("D.png", 29, "R", 0, 2, 1),
    @("J.png", 30, "R", 0, 3, 1),
    @("H.png", 31, "R", 0, 1, 1),
    @("B.png", 32, "R", 0, 2, 1),
    @("C.png", 33, "R", 0, 3, 1),
    @("M.png", 34, "R", 0, 1, 1),
    @("K.png", 35, "R", 0, 2, 1),
    @("L.png", 36, "R", 0, 3, 1),
    @("E.png", 37, "R", 0, 1, 1),
    @("F.png", 38, "R", 0, 2, 1),
    @("A.png", 39, "R", 0, 3, 1),
    @("L.png", 40, "R", 0, 1, 1),
    @("E.png", 41, "R", 0, 2, 1),
    @("H.png", 42, "R", 0, 3, 1),
    @("C.png", 43, "R", 0, 1, 1),
    @("G.png", 44, "R", 0, 2, 1),
    @("M.png", 45, "R", 0, 3, 1),
    @("H.png", 46, "R", 0, 1, 1),
    @("D.png", 47, "R", 0, 2, 1),
    @("E.png", 48, "R", 0, 3, 1)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# Match the saved selection state recorded in the diff
$ws.Range("A26:XFD49").Select() | Out-Null
